$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DailyData")

# --- November7 block (columns H:L) ---
# New "Matches" / "PPM" mini stat, two rows below the existing Averages table (H10:K11)
$ws.Range("H13").Value = "Matches"
$ws.Range("I13").Value = "PPM"
$ws.Range("H14").Formula = "=ROWS(November7[Points])"
$ws.Range("I14").Formula = "=November7[[#Totals],[Points]]/ROWS(November7[Points])"

# --- November10 block (columns V:Z) ---
$ws.Range("V14").Value = "Matches"
$ws.Range("W14").Value = "PPM"
$ws.Range("V15").Formula = "=ROWS(November10[Points])"
$ws.Range("W15").Formula = "=November10[[#Totals],[Points]]/ROWS(November10[Points])"

# --- November9 block (columns O:S) ---
$ws.Range("O17").Value = "Matches"
$ws.Range("P17").Value = "PPM"
$ws.Range("O18").Formula = "=ROWS(November9[Points])"
$ws.Range("P18").Formula = "=November9[[#Totals],[Points]]/ROWS(November9[Points])"

# --- November6 block (columns A:E) ---
$ws.Range("A29").Value = "Matches"
$ws.Range("B29").Value = "PPM"
$ws.Range("A30").Formula = "=ROWS(November6[Points])"
$ws.Range("B30").Formula = "=November6[[#Totals],[Points]]/ROWS(November6[Points])"

# Update selection/view to match the target file
$ws.Activate()
$ws.Application.ActiveWindow.ScrollColumn = 8
$ws.Range("W19").Select()

$ws1 = $wb.Worksheets.Item("Master")
$ws1.Range("G8").Select()
$ws.Activate()
